$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "0    (False, 14)`n1    (False, 30)`n2     (True, 12)`n3    (False, 16)`n4    (False, 12)`n5     (True, 18)`n6     (True, 13)`n7    (False, 12)`n8     (True, 12)`n9     (True, 27)`ndtype: object"
$ws.Range("D3").Value = "0    (False, 29)`n1     (True, 16)`n2     (True, 21)`n3     (True, 26)`n4     (True, 17)`n5    (False, 16)`n6    (False, 17)`n7    (False, 15)`n8     (True, 24)`n9    (False, 12)`ndtype: object"
$ws.Range("D4").Value = "0    (False, 26)`n1    (False, 18)`n2     (True, 30)`n3    (False, 23)`n4    (False, 30)`n5    (False, 15)`n6    (False, 27)`n7     (True, 30)`n8    (False, 28)`n9     (True, 19)`ndtype: object"
$ws.Range("D5").Value = "0     (True, 13)`n1    (False, 26)`n2     (True, 18)`n3     (True, 10)`n4     (True, 27)`n5    (False, 16)`n6     (True, 15)`n7    (False, 21)`n8     (True, 14)`n9    (False, 25)`ndtype: object"
$ws.Range("D6").Value = "0     (True, 18)`n1     (True, 23)`n2     (True, 25)`n3    (False, 23)`n4     (True, 25)`n5     (True, 29)`n6    (False, 15)`n7     (True, 30)`n8     (True, 25)`n9     (True, 27)`ndtype: object"
$ws.Range("D7").Value = "0     (True, 17)`n1     (True, 27)`n2    (False, 12)`n3     (True, 25)`n4    (False, 18)`n5    (False, 12)`n6    (False, 23)`n7     (True, 17)`n8     (True, 28)`n9    (False, 22)`ndtype: object"
$ws.Range("D8").Value = "0     (True, 28)`n1    (False, 29)`n2    (False, 21)`n3    (False, 23)`n4     (True, 14)`n5    (False, 13)`n6    (False, 11)`n7     (True, 19)`n8    (False, 25)`n9    (False, 18)`ndtype: object"
$ws.Range("D9").Value = "0     (True, 20)`n1     (True, 11)`n2    (False, 12)`n3    (False, 12)`n4    (False, 13)`n5    (False, 20)`n6    (False, 13)`n7     (True, 11)`n8     (True, 18)`n9    (False, 29)`ndtype: object"
$ws.Range("D10").Value = "0     (True, 15)`n1     (True, 11)`n2     (True, 10)`n3     (True, 12)`n4    (False, 14)`n5    (False, 28)`n6    (False, 29)`n7    (False, 19)`n8    (False, 11)`n9    (False, 15)`ndtype: object"
$ws.Range("D11").Value = "0    (False, 29)`n1    (False, 12)`n2     (True, 18)`n3    (False, 19)`n4    (False, 25)`n5    (False, 12)`n6    (False, 21)`n7     (True, 17)`n8    (False, 13)`n9     (True, 14)`ndtype: object"
$ws.Range("D12").Value = "0    (False, 26)`n1    (False, 25)`n2     (True, 28)`n3    (False, 24)`n4    (False, 16)`n5    (False, 29)`n6     (True, 17)`n7    (False, 14)`n8     (True, 12)`n9    (False, 15)`ndtype: object"
$ws.Range("D13").Value = "0    (False, 29)`n1    (False, 15)`n2    (False, 21)`n3    (False, 20)`n4    (False, 14)`n5    (False, 25)`n6    (False, 13)`n7     (True, 16)`n8    (False, 22)`n9    (False, 15)`ndtype: object"
$ws.Range("D14").Value = "0    (False, 28)`n1    (False, 27)`n2     (True, 25)`n3     (True, 28)`n4    (False, 24)`n5     (True, 29)`n6    (False, 29)`n7     (True, 20)`n8     (True, 28)`n9     (True, 17)`ndtype: object"
$ws.Range("D15").Value = "0     (True, 30)`n1     (True, 20)`n2     (True, 27)`n3     (True, 20)`n4     (True, 16)`n5    (False, 30)`n6    (False, 24)`n7     (True, 25)`n8    (False, 28)`n9     (True, 17)`ndtype: object"
$ws.Range("D16").Value = "0     (True, 17)`n1     (True, 11)`n2    (False, 12)`n3     (True, 16)`n4     (True, 20)`n5     (True, 29)`n6    (False, 20)`n7     (True, 11)`n8    (False, 18)`n9     (True, 28)`ndtype: object"
$ws.Range("D17").Value = "0     (True, 12)`n1    (False, 28)`n2    (False, 19)`n3     (True, 13)`n4    (False, 14)`n5    (False, 26)`n6    (False, 22)`n7     (True, 22)`n8    (False, 24)`n9    (False, 28)`ndtype: object"
$ws.Range("D18").Value = "0    (False, 23)`n1    (False, 27)`n2     (True, 10)`n3    (False, 17)`n4    (False, 18)`n5     (True, 15)`n6    (False, 20)`n7    (False, 21)`n8    (False, 30)`n9    (False, 10)`ndtype: object"
$ws.Range("D19").Value = "0     (True, 27)`n1     (True, 21)`n2    (False, 27)`n3     (True, 25)`n4    (False, 21)`n5    (False, 17)`n6     (True, 10)`n7    (False, 23)`n8    (False, 20)`n9     (True, 24)`ndtype: object"
$ws.Range("D20").Value = "0     (True, 23)`n1    (False, 14)`n2    (False, 24)`n3    (False, 17)`n4    (False, 12)`n5    (False, 15)`n6    (False, 16)`n7    (False, 24)`n8    (False, 17)`n9     (True, 22)`ndtype: object"
$ws.Range("D21").Value = "0    (False, 24)`n1     (True, 10)`n2     (True, 28)`n3     (True, 29)`n4     (True, 23)`n5     (True, 11)`n6     (True, 17)`n7    (False, 17)`n8    (False, 20)`n9     (True, 19)`ndtype: object"
$ws.Range("D22").Value = "0     (True, 25)`n1     (True, 21)`n2    (False, 17)`n3    (False, 12)`n4     (True, 21)`n5     (True, 28)`n6    (False, 23)`n7     (True, 14)`n8    (False, 22)`n9    (False, 11)`ndtype: object"
$ws.Range("D23").Value = "0     (True, 18)`n1    (False, 21)`n2     (True, 22)`n3     (True, 29)`n4    (False, 24)`n5     (True, 11)`n6     (True, 28)`n7    (False, 11)`n8     (True, 27)`n9     (True, 23)`ndtype: object"
$ws.Range("D24").Value = "0     (True, 12)`n1     (True, 14)`n2    (False, 13)`n3     (True, 16)`n4     (True, 28)`n5     (True, 14)`n6     (True, 13)`n7     (True, 28)`n8    (False, 13)`n9    (False, 29)`ndtype: object"
$ws.Range("D25").Value = "0     (True, 24)`n1     (True, 27)`n2     (True, 28)`n3     (True, 18)`n4     (True, 10)`n5    (False, 15)`n6     (True, 21)`n7    (False, 24)`n8    (False, 15)`n9     (True, 29)`ndtype: object"
$ws.Range("D26").Value = "0    (False, 22)`n1    (False, 28)`n2    (False, 14)`n3    (False, 29)`n4    (False, 11)`n5    (False, 16)`n6    (False, 19)`n7    (False, 25)`n8    (False, 13)`n9    (False, 12)`ndtype: object"
$ws.Range("D27").Value = "0     (True, 21)`n1    (False, 10)`n2     (True, 23)`n3     (True, 27)`n4     (True, 15)`n5     (True, 25)`n6    (False, 19)`n7     (True, 26)`n8    (False, 26)`n9    (False, 17)`ndtype: object"
$ws.Range("D28").Value = "0    (False, 28)`n1    (False, 11)`n2     (True, 13)`n3     (True, 26)`n4     (True, 12)`n5    (False, 17)`n6    (False, 29)`n7     (True, 18)`n8    (False, 15)`n9    (False, 17)`ndtype: object"
$ws.Range("D29").Value = "0    (False, 28)`n1    (False, 16)`n2    (False, 27)`n3    (False, 25)`n4     (True, 11)`n5    (False, 27)`n6    (False, 30)`n7     (True, 25)`n8     (True, 27)`n9     (True, 25)`ndtype: object"
$ws.Range("D30").Value = "0    (False, 17)`n1     (True, 28)`n2     (True, 14)`n3    (False, 18)`n4    (False, 21)`n5     (True, 27)`n6    (False, 22)`n7     (True, 19)`n8     (True, 18)`n9    (False, 11)`ndtype: object"
$ws.Range("D31").Value = "0     (True, 22)`n1    (False, 29)`n2     (True, 13)`n3     (True, 12)`n4    (False, 28)`n5     (True, 19)`n6    (False, 19)`n7     (True, 13)`n8    (False, 25)`n9    (False, 18)`ndtype: object"
$ws.Range("D32").Value = "0    (False, 21)`n1     (True, 23)`n2     (True, 16)`n3    (False, 24)`n4     (True, 11)`n5    (False, 18)`n6     (True, 30)`n7    (False, 22)`n8    (False, 30)`n9     (True, 10)`ndtype: object"
$ws.Range("D33").Value = "0    (False, 21)`n1    (False, 15)`n2     (True, 19)`n3     (True, 15)`n4    (False, 20)`n5     (True, 21)`n6    (False, 14)`n7     (True, 13)`n8    (False, 21)`n9    (False, 25)`ndtype: object"
$ws.Range("D34").Value = "0     (True, 27)`n1    (False, 25)`n2    (False, 27)`n3    (False, 27)`n4     (True, 29)`n5     (True, 10)`n6    (False, 16)`n7     (True, 23)`n8    (False, 11)`n9     (True, 28)`ndtype: object"
$ws.Range("D35").Value = "0     (True, 17)`n1    (False, 10)`n2     (True, 26)`n3     (True, 13)`n4     (True, 11)`n5     (True, 17)`n6     (True, 14)`n7    (False, 27)`n8     (True, 18)`n9    (False, 21)`ndtype: object"
$ws.Range("D36").Value = "0     (True, 28)`n1    (False, 15)`n2     (True, 13)`n3     (True, 29)`n4     (True, 11)`n5    (False, 12)`n6     (True, 27)`n7    (False, 18)`n8     (True, 27)`n9     (True, 14)`ndtype: object"
$ws.Range("D37").Value = "0    (False, 29)`n1    (False, 28)`n2     (True, 18)`n3    (False, 29)`n4     (True, 29)`n5    (False, 28)`n6    (False, 13)`n7    (False, 19)`n8     (True, 24)`n9    (False, 20)`ndtype: object"
$ws.Range("D38").Value = "0    (False, 13)`n1     (True, 19)`n2    (False, 10)`n3     (True, 12)`n4    (False, 28)`n5     (True, 19)`n6     (True, 19)`n7     (True, 12)`n8    (False, 15)`n9    (False, 15)`ndtype: object"
$ws.Range("D39").Value = "0    (False, 28)`n1     (True, 25)`n2     (True, 19)`n3    (False, 19)`n4    (False, 17)`n5     (True, 19)`n6    (False, 25)`n7    (False, 27)`n8    (False, 23)`n9     (True, 25)`ndtype: object"
$ws.Range("D40").Value = "0    (False, 10)`n1    (False, 14)`n2     (True, 12)`n3     (True, 20)`n4    (False, 19)`n5    (False, 15)`n6    (False, 11)`n7     (True, 27)`n8    (False, 27)`n9     (True, 18)`ndtype: object"
$ws.Range("D41").Value = "0    (False, 16)`n1    (False, 20)`n2    (False, 29)`n3    (False, 28)`n4    (False, 11)`n5     (True, 12)`n6     (True, 26)`n7     (True, 18)`n8     (True, 13)`n9     (True, 11)`ndtype: object"
